$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.974.25'
$ws.Range("E2").Value = '  -1.79%  '
$ws.Range("D3").Value = '1.820.06'
$ws.Range("E3").Value = '  -1.26%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.011'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.46%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.009'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4616'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.88%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3636'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07280'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8643'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.78'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.68%  '
$ws.Range("D12").Value = '1.882.90'
$ws.Range("E12").Value = '  +2.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07601'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.12'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.317'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.485'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.69%  '
$ws.Range("E17").Value = '  -0.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008616'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.009'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.48%  '
$ws.Range("D20").Value = '27.363.16'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.152'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.57'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.69%  '
$ws.Range("D24").Value = '2.111.56'
$ws.Range("E24").Value = '  +1.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.72'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.855'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.19'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.085'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.099'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.91'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08889'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.949'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.139'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7250'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.420'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.38%  '
$ws.Range("E36").Value = '  -0.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.491'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.95%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.074'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.93%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05272'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01915'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.928'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.148'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5189'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1631'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.249'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4843'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.81%  '
$ws.Range("E47").Value = '  -0.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.10'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.11'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.627'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.58%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06225'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.56%  '
